# Update column G ("K") values in the save_data sheet.
# These values are the computed "K" (strike count) values that were
# regenerated (regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 0
    6  = 0
    7  = 3
    8  = 1
    9  = 2
    10 = 1
    11 = 0
    12 = 1
    13 = 0
    14 = 1
    15 = 1
    16 = 2
    17 = 0
    18 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
